$d = $word.ActiveDocument

function New-PkgXml {
    param([string]$InnerParagraphXml)
    $pre = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>'
    $post = '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    return $pre + $InnerParagraphXml + $post
}

# 1) Remove the "Daniel Ojeda Velasco" heading run, keeping the now-empty
#    underlined paragraph mark in place (matches the commit: "Eliminado
#    cabecera de Reunion").
$p1 = $d.Paragraphs(1)
$r1 = $p1.Range
$r1.MoveEnd(1, -1) | Out-Null
$r1.Delete()

# 2) "-Se va a poder comer en el square, en la mitad(mamparas)" -> wrap
#    "square" with proofErr spellStart/spellEnd markers.
$frag3 = New-PkgXml '<w:r><w:t xml:space="preserve">-Se va a poder comer en el </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>square</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, en la mitad(mamparas)</w:t></w:r>'
$p3 = $d.Paragraphs(3)
$p3.Range.InsertXML($frag3)

# 3) "-Se van a poner enchufes en el square" -> wrap trailing "square"
#    with proofErr spellStart/spellEnd markers.
$frag4 = New-PkgXml '<w:r><w:t xml:space="preserve">-Se van a poner enchufes en el </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>square</w:t></w:r><w:proofErr w:type="spellEnd"/>'
$p4 = $d.Paragraphs(4)
$p4.Range.InsertXML($frag4)

# 4) "-Construir comunidad, sobretodo el primer año" -> wrap "sobretodo"
#    with proofErr spellStart/spellEnd markers.
$frag7 = New-PkgXml '<w:r><w:t xml:space="preserve">-Construir comunidad, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>sobretodo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> el primer año</w:t></w:r>'
$p7 = $d.Paragraphs(7)
$p7.Range.InsertXML($frag7)

# 5) "...Representación externa de la uni CEICEM..." -> wrap "uni" with
#    proofErr spellStart/spellEnd markers; keep the remaining runs as-is.
$frag12 = New-PkgXml '<w:r><w:t>-</w:t></w:r><w:r w:rsidR="00A93167"><w:t xml:space="preserve">Representación externa de la </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00A93167"><w:t>uni</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00356EEE"><w:t xml:space="preserve"> CEI</w:t></w:r><w:r w:rsidR="00E15C2F"><w:t>C</w:t></w:r><w:r w:rsidR="00356EEE"><w:t>E</w:t></w:r><w:r w:rsidR="00FA194B"><w:t xml:space="preserve">M </w:t></w:r><w:r w:rsidR="00C1139B"><w:t>(online, Madrid)</w:t></w:r><w:r w:rsidR="00356EEE"><w:t xml:space="preserve"> Y CEUNE</w:t></w:r><w:r w:rsidR="00FA194B"><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00C1139B"><w:t>(nacional)</w:t></w:r>'
$p12 = $d.Paragraphs(12)
$p12.Range.InsertXML($frag12)

Write-Output "done"
